$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source values (original D, J, K, L, M, P) keyed by original row number
$src = @{}
$src[2] = @{ D = 44676; J = 120; K = 4000; L = 4500; M = 4250; P = 71 }
$src[3] = @{ D = 44657; J = 100; K = 5000; L = 5500; M = 5250; P = 88 }
$src[4] = @{ D = 44648; J = 120; K = 6500; L = 7000; M = 6750; P = 112 }
$src[5] = @{ D = 44281; J = 120; K = 5500; L = 6000; M = 5750; P = 96 }
$src[6] = @{ D = 44935; J = 120; K = 6000; L = 7000; M = 6500; P = 108 }
$src[7] = @{ D = 44785; J = 130; K = 7000; L = 8000; M = 7500; P = 125 }
$src[8] = @{ D = 44362; J = 120; K = 8000; L = 9000; M = 8500; P = 142 }
$src[9] = @{ D = 44400; J = 120; K = 9000; L = 10000; M = 9500; P = 158 }
$src[10] = @{ D = 44421; J = 100; K = 8000; L = 9000; M = 8500; P = 142 }
$src[11] = @{ D = 44603; J = 140; K = 5500; L = 6000; M = 5750; P = 96 }
$src[12] = @{ D = 44740; J = 120; K = 6000; L = 7000; M = 6500; P = 108 }
$src[13] = @{ D = 44764; J = 120; K = 7000; L = 8000; M = 7500; P = 125 }
$src[14] = @{ D = 44669; J = 130; K = 4500; L = 5000; M = 4750; P = 79 }
$src[15] = @{ D = 44589; J = 110; K = 5000; L = 6000; M = 5500; P = 92 }
$src[16] = @{ D = 44242; J = 160; K = 5000; L = 5500; M = 5250; P = 88 }
$src[17] = @{ D = 44494; J = 120; K = 5000; L = 6000; M = 5500; P = 92 }
$src[18] = @{ D = 44760; J = 130; K = 7000; L = 7500; M = 7250; P = 121 }
$src[19] = @{ D = 44627; J = 120; K = 4000; L = 4500; M = 4250; P = 71 }
$src[20] = @{ D = 44382; J = 160; K = 7000; L = 8000; M = 7438; P = 124 }
$src[21] = @{ D = 44827; J = 120; K = 6000; L = 7000; M = 6500; P = 108 }

# Mapping: new row -> source row to pull values from
$mapping = @{
    2 = 12
    3 = 19
    4 = 20
    5 = 8
    6 = 11
    7 = 14
    8 = 7
    9 = 16
    10 = 2
    11 = 18
    12 = 15
    13 = 9
    14 = 4
    15 = 21
    16 = 10
    17 = 13
    18 = 17
    19 = 6
    20 = 5
    21 = 3
}

foreach ($row in $mapping.Keys) {
    $s = $src[$mapping[$row]]
    $ws.Range("D$row").Value = $s.D
    $ws.Range("J$row").Value = $s.J
    $ws.Range("K$row").Value = $s.K
    $ws.Range("L$row").Value = $s.L
    $ws.Range("M$row").Value = $s.M
    $ws.Range("P$row").Value = $s.P
}
